# Update the "Generate Report for Handback" timestamps in the handback-status workbook.
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for df692b68-...md moves forward ~56s
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 01:00:25"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the df692b68 -> zh-cn xlf row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 01:00:20"
$wsZhCn.Range("K2").Value = "2016-08-16 01:00:37"

# de-de sheet: Correspond Handoff datetime for the df692b68 -> de-de xlf row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 01:00:44"
